$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1441
$ws1.Range("F5").Value = 219
$ws1.Range("F6").Value = 88
$ws1.Range("F7").Value = 131
$ws1.Range("F8").Value = 6123
$ws1.Range("F12").Value = 4971
$ws1.Range("F22").Value = 3471

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 69

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 69
$ws4.Range("F5").Value = 1441
$ws4.Range("F6").Value = 219
$ws4.Range("F7").Value = 88
$ws4.Range("F8").Value = 131
$ws4.Range("F9").Value = 6123
$ws4.Range("F13").Value = 4971
$ws4.Range("F16").Value = 1165
$ws4.Range("F23").Value = 3471
